# Applies the cryptos-list refresh described by the diff: updates the
# Price (D) and Volume(1h) (E) values for the listed coin rows, and fixes
# the ordering of WrappedEther/Polkadot (rows 12-13) and RenderToken/
# Mantle (rows 47-48), whose rank positions swapped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.142.04"
$ws.Range("E2").Value = "  -0.51%  "

# Row 3
$ws.Range("D3").Value = "1.657.78"
$ws.Range("E3").Value = "  -0.59%  "

# Row 4
$ws.Range("E4").Value = "  -0.33%  "

# Row 5
$ws.Range("D5").Value = "'218.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.19%  "

# Row 6
$ws.Range("D6").Value = "'0.5286"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.99%  "

# Row 7
$ws.Range("E7").Value = "  -0.27%  "

# Row 8
$ws.Range("D8").Value = "'0.2606"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.24%  "

# Row 9
$ws.Range("D9").Value = "'0.06348"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.45%  "

# Row 10
$ws.Range("D10").Value = "'20.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.55%  "

# Row 11
$ws.Range("D11").Value = "'0.07790"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.45%  "

# Row 12
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'4.509"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.31%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.648.42"
$ws.Range("E13").Value = "  -1.33%  "

# Row 14
$ws.Range("D14").Value = "'0.5487"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.28%  "

# Row 15
$ws.Range("D15").Value = "0.0₅8218"
$ws.Range("E15").Value = "  -0.77%  "

# Row 16
$ws.Range("D16").Value = "'65.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.93%  "

# Row 17
$ws.Range("D17").Value = "26.148.39"
$ws.Range("E17").Value = "  -0.51%  "

# Row 18
$ws.Range("E18").Value = "  -0.35%  "

# Row 19
$ws.Range("D19").Value = "'4.585"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.89%  "

# Row 20
$ws.Range("D20").Value = "'193.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.70%  "

# Row 21
$ws.Range("D21").Value = "'10.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.41%  "

# Row 22
$ws.Range("D22").Value = "'6.042"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.46%  "

# Row 23
$ws.Range("D23").Value = "'1.003"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.44%  "

# Row 24
$ws.Range("D24").Value = "'141.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.87%  "

# Row 25
$ws.Range("D25").Value = "'0.1254"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.36%  "

# Row 26
$ws.Range("D26").Value = "'7.285"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.22%  "

# Row 27
$ws.Range("D27").Value = "'16.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.05%  "

# Row 29
$ws.Range("D29").Value = "'0.05954"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.68%  "

# Row 30
$ws.Range("D30").Value = "'1.281"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.01%  "

# Row 31
$ws.Range("D31").Value = "'3.526"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.81%  "

# Row 32
$ws.Range("D32").Value = "'3.263"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.92%  "

# Row 33
$ws.Range("D33").Value = "'1.580"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.11%  "

# Row 34
$ws.Range("D34").Value = "'0.9542"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.74%  "

# Row 35
$ws.Range("D35").Value = "'2.795"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.47%  "

# Row 36
$ws.Range("E36").Value = "  -0.64%  "

# Row 37
$ws.Range("D37").Value = "'0.5689"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.95%  "

# Row 38
$ws.Range("D38").Value = "'0.01616"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.31%  "

# Row 39
$ws.Range("D39").Value = "'5.827"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.13%  "

# Row 40
$ws.Range("D40").Value = "'0.8500"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.89%  "

# Row 42
$ws.Range("D42").Value = "'103.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.93%  "

# Row 43
$ws.Range("D43").Value = "1.026.66"
$ws.Range("E43").Value = "  +0.51%  "

# Row 44
$ws.Range("D44").Value = "1.801.10"
$ws.Range("E44").Value = "  -0.46%  "

# Row 45
$ws.Range("D45").Value = "'57.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.36%  "

# Row 46
$ws.Range("D46").Value = "'1.010"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.37%  "

# Row 47
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "'1.499"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.10%  "

# Row 48
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.4290"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.67%  "

# Row 49
$ws.Range("D49").Value = "'0.05153"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.62%  "

# Row 50
$ws.Range("D50").Value = "'7.807"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.47%  "

# Row 51
$ws.Range("D51").Value = "'0.09715"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.47%  "
